$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 79
$ws.Range("F3").Value = 211
$ws.Range("F4").Value = 98
$ws.Range("F6").Value = 3285
$ws.Range("F7").Value = 908
$ws.Range("F8").Value = 2103
$ws.Range("F9").Value = 2020
$ws.Range("F10").Value = 1046
$ws.Range("F11").Value = 539
$ws.Range("F13").Value = 1631
$ws.Range("F14").Value = 352
$ws.Range("F16").Value = 18
$ws.Range("F17").Value = 77
$ws.Range("F18").Value = 105
$ws.Range("F19").Value = 1473
$ws.Range("F20").Value = 553
$ws.Range("F21").Value = 656
$ws.Range("F22").Value = 340
$ws.Range("F23").Value = 11867
$ws.Range("F24").Value = 11891
$ws.Range("F25").Value = 872
$ws.Range("G25").Value = 54
$ws.Range("G26").Value = 49.5
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 54
$ws.Range("F28").Value = 1866
$ws.Range("F29").Value = 166
$ws.Range("F30").Value = 479

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 37

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 79
$ws.Range("F4").Value = 211
$ws.Range("F6").Value = 98
$ws.Range("F8").Value = 3285
$ws.Range("F9").Value = 908
$ws.Range("F10").Value = 2103
$ws.Range("F11").Value = 2020
$ws.Range("F12").Value = 1046
$ws.Range("F13").Value = 539
$ws.Range("F15").Value = 1631
$ws.Range("F16").Value = 352
$ws.Range("F18").Value = 18
$ws.Range("F20").Value = 77
$ws.Range("F21").Value = 37
$ws.Range("F22").Value = 105
$ws.Range("F23").Value = 1473
$ws.Range("F24").Value = 553
$ws.Range("F25").Value = 656
$ws.Range("F26").Value = 340
$ws.Range("F27").Value = 11867
$ws.Range("F28").Value = 11892
$ws.Range("F29").Value = 872
$ws.Range("G29").Value = 54
$ws.Range("G30").Value = 49.5
$ws.Range("F31").Value = 6
$ws.Range("G31").Value = 54
$ws.Range("F32").Value = 1866
$ws.Range("F35").Value = 166
$ws.Range("F36").Value = 479
